$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report header text in G1 from "...for the month" to "...for the period"
$ws.Range("G1").Value = "Total fees accrued for the period"

# Reflect the saved selection state (active cell moved to G2)
$ws.Range("G2").Select()
